# Actualización automática 2025-06-11 11:20:08
# Remove the client "JACOME MONCAYO XAVIER ALFONSO" (row 29) from both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting subsequent rows
# up by one, and refresh the trailing totals row accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO -------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(29).Delete()

# The trailing summary row (now row 55) holds text like "1 de 54"; the
# count of clients dropped from 54 to 53, so update the denominator while
# keeping the numerator unchanged.
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N")
foreach ($col in $cols1) {
    $cell = $ws1.Range($col + "55")
    $text = $cell.Value2
    $cell.Value = $text.Replace("de 54", "de 53")
}

# --- Sheet: VENTA MENSUAL ------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(29).Delete()

# The trailing totals row (now row 55) holds numeric sums that must be
# reduced by the removed client's figures.
$ws2.Range("C55").Value = 87328.29000000001
$ws2.Range("D55").Value = 94185.91
$ws2.Range("E55").Value = 54896.95
$ws2.Range("F55").Value = 29609.21
$ws2.Range("G55").Value = 88000
